# Update the daily scores sheet:
#  - correct the 2025-02-23 abs_activity score (row 90, cols C & F)
#  - append two new days of data (2025-02-24 and 2025-02-25), four rows each
#    (abs_activity, rel_activity, abs_sleep, rel_sleep), in rows 94-101

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 90 (2025-02-23 / abs_activity) ---
$ws.Cells.Item(90, 3).Value = 9.224240567669321
$ws.Cells.Item(90, 6).Value = 9.224240567669321

# --- New data rows ---
$newRows = @(
    @{ Row = 94;  Date = "2025-02-24"; Typ = "abs_activity"; C = 10;                D = 0; E = 0; F = 10 },
    @{ Row = 95;  Date = "2025-02-24"; Typ = "rel_activity"; C = 0;                 D = 0; E = 0; F = 0 },
    @{ Row = 96;  Date = "2025-02-24"; Typ = "abs_sleep";    C = 4.133333333333333; D = 0; E = 0; F = 4.133333333333333 },
    @{ Row = 97;  Date = "2025-02-24"; Typ = "rel_sleep";    C = 0;                 D = 0; E = 0; F = 0 },
    @{ Row = 98;  Date = "2025-02-25"; Typ = "abs_activity"; C = 7.062117055431884; D = 0; E = 0; F = 7.062117055431884 },
    @{ Row = 99;  Date = "2025-02-25"; Typ = "rel_activity"; C = 0;                 D = 0; E = 0; F = 0 },
    @{ Row = 100; Date = "2025-02-25"; Typ = "abs_sleep";    C = 9.633333333333333; D = 0; E = 0; F = 9.633333333333333 },
    @{ Row = 101; Date = "2025-02-25"; Typ = "rel_sleep";    C = 0;                 D = 0; E = 0; F = 0 }
)

foreach ($r in $newRows) {
    # The Date column holds plain text like "2025-02-24" in this workbook (not a
    # real Excel date). Assigning that string straight to .Value would make Excel
    # auto-detect it as a date and stamp a date number format on the cell, so we
    # pre-format the cell as Text, write the value, then clear the formatting back
    # to the default (unstyled) cell - this keeps the stored type as text without
    # leaving a residual number format behind.
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.Typ
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
